$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.567.05'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.919.02'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.67'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E7').Value = '  +3.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2914'
$ws.Range('E8').Value = '  +1.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06728'
$ws.Range('E9').Value = '  -2.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '107.26'
$ws.Range('E10').Value = '  +0.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '18.59'
$ws.Range('E11').Value = '  +1.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.920.42'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07670'
$ws.Range('E13').Value = '  +0.41%  '
$ws.Range('E14').Value = '  +2.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6676'
$ws.Range('E15').Value = '  +1.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '281.12'
$ws.Range('E16').Value = '  -5.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.544.79'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007557'
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.170.64'
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.81'
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.497'
$ws.Range('E22').Value = '  +4.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.455'
$ws.Range('E24').Value = '  +3.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.472'
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.64'
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.14'
$ws.Range('E27').Value = '  -6.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.116'
$ws.Range('E28').Value = '  +3.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1059'
$ws.Range('E29').Value = '  -3.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.406'
$ws.Range('E30').Value = '  +3.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.156'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.053'
$ws.Range('E32').Value = '  +1.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05016'
$ws.Range('E33').Value = '  -1.54%  '
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.140'
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9998'
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02037'
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.683'
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '111.64'
$ws.Range('E40').Value = '  +3.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.019'
$ws.Range('E41').Value = '  -1.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4458'
$ws.Range('E42').Value = '  +5.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8724'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.927'
$ws.Range('E44').Value = '  +1.79%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '68.16'
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '49.40'
$ws.Range('E47').Value = '  -7.11%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.272'
$ws.Range('E48').Value = '  +0.92%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.336'
$ws.Range('E49').Value = '  +1.13%  '
$ws.Range('E50').Value = '  +3.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.93'
$ws.Range('E51').Value = '  +0.74%  '
